$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44964
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 20000
$ws.Range("L2").Value = 21000
$ws.Range("M2").Value = 20500
$ws.Range("P2").Value = 1139

# Row 3
$ws.Range("D3").Value = 45229
$ws.Range("J3").Value = 460
$ws.Range("K3").Value = 16000
$ws.Range("L3").Value = 17000
$ws.Range("M3").Value = 16500
$ws.Range("P3").Value = 917

# Row 4
$ws.Range("D4").Value = 44984
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 17000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 17500
$ws.Range("P4").Value = 972

# Row 5
$ws.Range("D5").Value = 44957
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 21000
$ws.Range("L5").Value = 22000
$ws.Range("M5").Value = 21500
$ws.Range("P5").Value = 1194

# Row 6
$ws.Range("D6").Value = 45177
$ws.Range("J6").Value = 540

# Row 7
$ws.Range("D7").Value = 45117
$ws.Range("J7").Value = 300

# Row 8
$ws.Range("D8").Value = 45180
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 16500
$ws.Range("L8").Value = 17000
$ws.Range("M8").Value = 16750
$ws.Range("P8").Value = 931

# Row 9
$ws.Range("D9").Value = 45154
$ws.Range("J9").Value = 500

# Row 10
$ws.Range("D10").Value = 45222
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 16000
$ws.Range("L10").Value = 17000
$ws.Range("M10").Value = 16500
$ws.Range("P10").Value = 917

# Row 11
$ws.Range("D11").Value = 44568
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 16000
$ws.Range("M11").Value = 15500
$ws.Range("P11").Value = 861

# Row 12
$ws.Range("D12").Value = 44557
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 13000
$ws.Range("L12").Value = 14000
$ws.Range("M12").Value = 13500
$ws.Range("P12").Value = 750

# Row 13
$ws.Range("D13").Value = 44960
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 19500
$ws.Range("L13").Value = 20000
$ws.Range("M13").Value = 19750
$ws.Range("P13").Value = 1097

# Row 14
$ws.Range("D14").Value = 45215
$ws.Range("J14").Value = 400
$ws.Range("K14").Value = 16000
$ws.Range("L14").Value = 17000
$ws.Range("M14").Value = 16500
$ws.Range("P14").Value = 917

# Row 15
$ws.Range("D15").Value = 44977
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 16500
$ws.Range("M15").Value = 16750
$ws.Range("P15").Value = 931

# Row 16
$ws.Range("D16").Value = 44998
$ws.Range("J16").Value = 320
$ws.Range("K16").Value = 17000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 17500
$ws.Range("P16").Value = 972

# Row 17
$ws.Range("D17").Value = 45194
$ws.Range("J17").Value = 400
$ws.Range("K17").Value = 16500
$ws.Range("M17").Value = 16750
$ws.Range("P17").Value = 931

# Row 18
$ws.Range("D18").Value = 44547
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 13000
$ws.Range("L18").Value = 14000
$ws.Range("M18").Value = 13500
$ws.Range("P18").Value = 750

# Row 19
$ws.Range("D19").Value = 45159
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 16000
$ws.Range("L19").Value = 17000
$ws.Range("M19").Value = 16500
$ws.Range("P19").Value = 917

# Row 20
$ws.Range("D20").Value = 45230
$ws.Range("J20").Value = 360
$ws.Range("K20").Value = 16000
$ws.Range("L20").Value = 17000
$ws.Range("M20").Value = 16500
$ws.Range("P20").Value = 917

# Row 21
$ws.Range("D21").Value = 45142
$ws.Range("J21").Value = 400
$ws.Range("K21").Value = 17000
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = 17500
$ws.Range("P21").Value = 972

# Row 22
$ws.Range("D22").Value = 45068

# Row 24
$ws.Range("D24").Value = 45152
$ws.Range("J24").Value = 500

# Row 25
$ws.Range("D25").Value = 45005
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = 17000
$ws.Range("L25").Value = 18000
$ws.Range("M25").Value = 17500
$ws.Range("P25").Value = 972
